# Update factsheets with text edits from COMM
#
# The "count" values on each factsheet tab were being stored as numbers;
# convert them to literal text instead (so they round-trip exactly, e.g.
# through templating), and append a "Total" row to the County tab that
# mirrors the Overall tab's totals.
#
# NumberFormat is forced to "@" (Text) right before the assignment so the
# numeric-looking strings ("547", "23", ...) aren't silently re-parsed back
# into numbers by Excel, then the style is reset to Normal afterwards so no
# stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ----- Overall sheet -----
$ws = $wb.Worksheets.Item("Overall")
Set-TextValue $ws.Range("A2") "547"

# ----- County sheet -----
$ws = $wb.Worksheets.Item("County")
Set-TextValue $ws.Range("B2") "23"
Set-TextValue $ws.Range("B3") "54"
Set-TextValue $ws.Range("B4") "59"
Set-TextValue $ws.Range("B5") "343"
Set-TextValue $ws.Range("B6") "68"

Set-TextValue $ws.Range("A7") "Total"
Set-TextValue $ws.Range("B7") "547"
Set-TextValue $ws.Range("C7") '$1,009,094,930'
Set-TextValue $ws.Range("D7") "9.60%"
Set-TextValue $ws.Range("E7") "-11.45%"
Set-TextValue $ws.Range("F7") "66.00%"

# ----- Congressional District sheet -----
$ws = $wb.Worksheets.Item("Congressional District")
Set-TextValue $ws.Range("B2") "281"
Set-TextValue $ws.Range("B3") "266"
Set-TextValue $ws.Range("B4") "547"

# ----- Size sheet -----
$ws = $wb.Worksheets.Item("Size")
Set-TextValue $ws.Range("B2") "177"
Set-TextValue $ws.Range("B3") "155"
Set-TextValue $ws.Range("B4") "79"
Set-TextValue $ws.Range("B5") "45"
Set-TextValue $ws.Range("B6") "64"
Set-TextValue $ws.Range("B7") "27"
Set-TextValue $ws.Range("B8") "547"

# ----- Subsector sheet -----
$ws = $wb.Worksheets.Item("Subsector")
Set-TextValue $ws.Range("B2") "63"
Set-TextValue $ws.Range("B3") "82"
Set-TextValue $ws.Range("B4") "26"
Set-TextValue $ws.Range("B5") "42"
Set-TextValue $ws.Range("B6") "2"
Set-TextValue $ws.Range("B7") "153"
Set-TextValue $ws.Range("B8") "5"
Set-TextValue $ws.Range("B9") "2"
Set-TextValue $ws.Range("B10") "34"
Set-TextValue $ws.Range("B11") "3"
Set-TextValue $ws.Range("B12") "129"
Set-TextValue $ws.Range("B13") "6"
Set-TextValue $ws.Range("B14") "547"
